# Fruta / hortaliza, semanal
# Insert a new weekly record at row 150 (shifting existing rows 150:181 down to
# 151:182), then populate the new row with this week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 150 (copy + insert) so the new row keeps identical formatting
# (including the date style on column D) and all of the fixed/common column
# values (Mercado, Region, Tipo, Producto, etc.) that are shared by every
# record in this block. This also naturally shifts rows 150:181 down to
# 151:182.
$ws.Rows("150:150").Copy()
$ws.Rows("150:150").Insert()

# Now overwrite the cells that differ for the new weekly record in row 150.
$ws.Range("D150").Value2 = 44644
$ws.Range("L150").Value2 = "Primera"
$ws.Range("M150").Value2 = 200
$ws.Range("N150").Value2 = 8000
$ws.Range("O150").Value2 = 8500
$ws.Range("P150").Value2 = 8250
$ws.Range("R150").Value2 = "Perú"
$ws.Range("S150").Value2 = 2062
$ws.Range("T150").Value2 = 4
